$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before current row 4 (shifts rows 4-6 and 9-11 down by 1)
$ws.Rows.Item(4).Insert()

# Match formatting of the new blank row4 cells to row3 (same blank-row style as before)
$ws.Range("A3:B3").Copy()
$ws.Range("A4:B4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Add "Done" markers in column D for rows 2 and 3
$ws.Range("D2").Value = "Done"
$ws.Range("D3").Value = "Done"

# Update selection to match target
$ws.Range("C4").Select()
